$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Price Target | Recs")

# Insert a new row at 23 (pushes the old rows 23-40 down to 24-41),
# matching the existing row's formatting (dates/number/percent styles,
# shared formulas, etc. all shift automatically).
$ws.Rows("23:23").Insert()

# Fill in the newly inserted row 23 with the NVDA trade.
$ws.Range("B23").Value = "NVDA"
$ws.Range("C23").Value = 45754
$ws.Range("D23").Value = "Long"
$ws.Range("E23").Value = 96.69
$ws.Range("F23").Value = 7

# RNMBY (row 34 after the shift): status changed from Long to Market Outpreform.
$ws.Range("D34").Value = "Market Outpreform"

# PEN (row 35 after the shift): trade closed out with an exit price/return/note.
$ws.Range("G35").Value = 55.21
$ws.Range("H35").Formula = "=G35/E35-1"
$ws.Range("I35").Value = "Quick in and out tactictal trade … need to commit to further research, but I like paypal more from a value perspective"
$ws.Range("J34").Copy()
$ws.Range("J35").PasteSpecial(-4122)
$ws.Range("J35").Value = "Y"

# SRPT (row 37 after the shift): conviction/size updated from 4 to 8.
$ws.Range("F37").Value = 8

# New trades appended below the existing data.
$ws.Range("C23").Copy()

$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("B38").Value = "AMD"
$ws.Range("C38").Value = 45790
$ws.Range("D38").Value = "Long"
$ws.Range("E38").Value = 116.27
$ws.Range("F38").Value = 4

$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("B39").Value = "EL"
$ws.Range("C39").Value = 45791
$ws.Range("D39").Value = "Long"
$ws.Range("E39").Value = 62.11
$ws.Range("F39").Value = 4

$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("B40").Value = "LULU"
$ws.Range("C40").Value = 45791
$ws.Range("D40").Value = "Long"
$ws.Range("E40").Value = 313.86
$ws.Range("F40").Value = 5

$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("B41").Value = "QBTS"
$ws.Range("C41").Value = 45796
$ws.Range("D41").Value = "Short"
$ws.Range("E41").Value = 13.24
$ws.Range("F41").Value = 6

# Restore the view state recorded in the saved workbook.
$ws.Range("A28").Select()
$ws.Range("G41").Select()

$wb.Windows.Item(1).WindowState = -4143
